$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 6264, 46069.95833333334),
    @(3, 6226, 46069.96875),
    @(4, 6217, 46069.97916666666),
    @(5, 6139, 46069.98958333334),
    @(6, 6094, 46070),
    @(7, 6046, 46070.01041666666),
    @(8, 6042, 46070.03125),
    @(9, 6008, 46070.04166666666),
    @(10, 5951, 46070.05208333334),
    @(11, 5938, 46070.0625),
    @(12, 5920, 46070.07291666666),
    @(13, 5945, 46070.08333333334),
    @(14, 5967, 46070.09375),
    @(15, 5949, 46070.10416666666),
    @(16, 5917, 46070.11458333334),
    @(17, 5910, 46070.125),
    @(18, 6021, 46070.13541666666),
    @(19, 6059, 46070.14583333334),
    @(20, 6056, 46070.15625),
    @(21, 6100, 46070.16666666666),
    @(22, 6221, 46070.17708333334),
    @(23, 6391, 46070.1875),
    @(24, 6459, 46070.19791666666),
    @(25, 6727, 46070.20833333334),
    @(26, 6880, 46070.21875),
    @(27, 6975, 46070.22916666666),
    @(28, 7071, 46070.23958333334),
    @(29, 7360, 46070.25),
    @(30, 7480, 46070.26041666666),
    @(31, 7661, 46070.27083333334),
    @(32, 7804, 46070.28125),
    @(33, 8068, 46070.29166666666)
)

foreach ($row in $data) {
    $ws.Cells.Item($row[0], 1).Value = $row[1]
    $ws.Cells.Item($row[0], 2).Value = $row[2]
}
